$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '28.663.37'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.26%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.871.82'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.25%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.006'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.49%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.82'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.01%  '

$ws.Range('E6').Value = '  +0.47%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4670'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.11%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3888'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.50%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07874'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.09%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9738'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.57%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '22.00'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.48%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.842.07'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -5.00%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.998'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.59%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.710'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.92%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06985'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.30%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.18'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.65%  '

$ws.Range('E17').Value = '  +0.66%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001005'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.16%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.82'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.04%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.005'
$ws.Range('D20').Style = 'Normal'

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '28.659.12'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.11%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.302'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.05%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.02'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.17%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.117'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.78%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.077.89'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.21%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '152.91'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.46%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.24'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.25%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.745'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.24%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.36'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.77%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09372'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.38%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9204'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.78%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.274'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.38%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.340'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.87%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.344'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.71%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05815'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.91%  '

$ws.Range('E37').Value = '  -1.96%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.148'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.33%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.762'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.25%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5628'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.84%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1788'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.60%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.764'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.28%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.07223'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.76%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.71'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.24%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5314'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.96%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.156'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.90%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.825'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.39%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '113.17'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.22%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.064'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.96%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.414'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.06%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.006'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.60%  '
